$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.461
$ws.Range("C6").Value = -12.974
$ws.Range("C7").Value = -13.057
$ws.Range("B8").Value = 5.314000000000001
$ws.Range("C8").Value = -12.154
$ws.Range("A12").Value = -21.531
$ws.Range("B12").Value = 6.842000000000001
$ws.Range("B14").Value = 6.622
$ws.Range("C19").Value = -13.057
$ws.Range("C21").Value = -13.036
$ws.Range("B22").Value = 6.74
$ws.Range("C24").Value = -12.512
